$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix row 4 (Trạm y tế Mỹ Long): Mã trạm y tế = 1
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = "Trạm y tế Mỹ Long"

# Fix row 5: now the first of the newly added "Phường" health stations
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = "Trạm y tế Phường Mỹ Bình"

# New rows 6-16: additional ward/commune health stations
$newStations = @(
    @{ Row = 6;  D = 3;  Name = "Trạm y tế Phường Mỹ Xuyên" },
    @{ Row = 7;  D = 4;  Name = "Trạm y tế Phường Đông Xuyên" },
    @{ Row = 8;  D = 5;  Name = "Trạm y tế Mỹ Hòa" },
    @{ Row = 9;  D = 6;  Name = "Trạm y tế Phường Mỹ Phước" },
    @{ Row = 10; D = 7;  Name = "Trạm y tế Phường Mỹ Quý" },
    @{ Row = 11; D = 8;  Name = "Trạm y tế Phường Mỹ Thới" },
    @{ Row = 12; D = 9;  Name = "Trạm y tế Phường Mỹ Thạnh" },
    @{ Row = 13; D = 10; Name = "Trạm y tế Phường Bình Khánh" },
    @{ Row = 14; D = 11; Name = "Trạm y tế Phường Bình Đức" },
    @{ Row = 15; D = 12; Name = "Trạm y tế Xã Mỹ Khánh" },
    @{ Row = 16; D = 13; Name = "Trạm y tế Xã Mỹ Hòa Hưng" }
)

foreach ($station in $newStations) {
    $r = $station.Row
    if ($r -lt 16) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
    $ws.Cells.Item($r, 2).Value = 5
    $ws.Cells.Item($r, 3).Value = 78
    $ws.Cells.Item($r, 4).Value = $station.D
    $ws.Cells.Item($r, 5).Value = $station.Name
}

$ws.Cells.Item(17, 3).Select()
